# Atualizado por script em 12-11-2023 20:45
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Swap the match data (columns F:V) between the following row
#    pairs. Columns A:E (Indice/pais/torneio/temporada/data_partida)
#    stay attached to their original row number.
# -----------------------------------------------------------------
$pairs = @(
    @(8, 9),
    @(16, 17),
    @(22, 23)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = $ws.Range("F$r1`:V$r1")
    $range2 = $ws.Range("F$r2`:V$r2")
    $v1 = $range1.Value2
    $v2 = $range2.Value2
    $range1.Value2 = $v2
    $range2.Value2 = $v1
}

# -----------------------------------------------------------------
# 2) Append a new row (48) with a new match record, copying the
#    formatting (styles) from row 47 so number formats / fonts /
#    borders match the rest of the table.
# -----------------------------------------------------------------
$ws.Range("A47:V47").Copy($ws.Range("A48:V48"))

$ws.Range("A48").Value2 = 47
$ws.Range("B48").Value2 = "lebanon"
$ws.Range("C48").Value2 = "premier-league"
$ws.Range("D48").Value2 = "2023-2024"
$ws.Range("E48").Value2 = 45242.65625
$ws.Range("F48").Value2 = "Racing"
$ws.Range("G48").Value2 = 1
$ws.Range("H48").Value2 = "Al Sahel"
$ws.Range("I48").Value2 = 2
$ws.Range("J48").Value2 = 2.42
$ws.Range("K48").Value2 = "11/11/2023 04:12"
$ws.Range("L48").Value2 = 2.81
$ws.Range("M48").Value2 = "12/11/2023 14:37"
$ws.Range("N48").Value2 = 3.09
$ws.Range("O48").Value2 = "11/11/2023 04:12"
$ws.Range("P48").Value2 = 3.15
$ws.Range("Q48").Value2 = "12/11/2023 13:50"
$ws.Range("R48").Value2 = 2.65
$ws.Range("S48").Value2 = "11/11/2023 04:12"
$ws.Range("T48").Value2 = 2.52
$ws.Range("U48").Value2 = "12/11/2023 13:52"
$ws.Range("V48").Value2 = "https://www.betexplorer.com/football/lebanon/premier-league/racing-al-sahel/WCQPFqWN/"
